# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP), G (sum)
# F (Win) is unchanged by this edit.
$data = @{
    2 = @{ B = 1.459612070389937;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732;  G = 3.781711156805759  }
    3 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 3.900430680208489;  E = 0.496779210170732;  G = 9.295990156953671  }
    4 = @{ B = 1.459612070389937;  C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732;  G = 4.429675500412797  }
    5 = @{ B = 1.459612070389937;  C = 1.667794583268128; D = 0.8054896365839992; E = 8.660232485948974;  G = 12.59312877619104  }
    6 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144  }
    7 = @{ B = 0.3048080303191223; C = 225321.0684179339; D = 0.1575252929769615; E = 8.660232485948974;  G = 225330.1909837431 }
    8 = @{ B = 1.459612070389937;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732;  G = 3.781711156805759  }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
